$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 97, shifting rows 97:125 down to 98:126
$ws.Rows.Item(97).Insert()

# Fill in the new row 97 data
$ws.Cells.Item(97, 1).Value = 5
$ws.Cells.Item(97, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(97, 3).Value = "Maule"
$ws.Cells.Item(97, 4).Value = 44468
$ws.Cells.Item(97, 5).Value = 7
$ws.Cells.Item(97, 6).Value = 100112021
$ws.Cells.Item(97, 7).Value = "Ají"
$ws.Cells.Item(97, 8).Value = "Americana (o)"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 100
$ws.Cells.Item(97, 11).Value = 80000
$ws.Cells.Item(97, 12).Value = 80000
$ws.Cells.Item(97, 13).Value = 80000
$ws.Cells.Item(97, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(97, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(97, 16).Value = 3200
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"
